$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "audiology"
$ws.Range("B2").Value = "https://www.ucl.ac.uk/prospective-students/graduate/taught-degrees/advanced-audiology-msc"
$ws.Range("C2").Value = "Advanced Audiology is designed for practising audiologists looking to enhance their clinical skills. Graduates are eligible to apply for senior clinical roles in the NHS or private sector. A variety of specialist modules can be selected to suit your professional needs. The programme has a strong research and evidence-based practice foundation."
$ws.Range("D2").Value = "English language requirements:https://www.ucl.ac.uk/prospective-students/graduate/learning-and-living-ucl/international-students/english-language-requirements"
$ws.Range("E2").Value = "UCL Graduate Admissions team:https://www.ucl.ac.uk/prospective-students/graduate/admissions-enquiries#form"
$ws.Range("F2").Value = "Application fees:https://www.ucl.ac.uk/prospective-students/graduate/application-fees"
$ws.Range("G2").Value = "Entry requirements:#entry-requirements"
$ws.Range("H2").Value = "Graduate degrees:/prospective-students/graduate/graduate-degrees"
$ws.Range("I2").Value = "Taught Degrees:/prospective-students/graduate/taught-degrees"
$ws.Range("J2").Value = "Applying for Graduate Taught Study at UCL:/prospective-students/graduate/taught-degrees/applying-graduate-taught-study-ucl"
$ws.Range("K2").Value = "Research Degrees:/prospective-students/graduate/research-degrees"
$ws.Range("L2").Value = "Applying for Graduate Research Study at UCL:/prospective-students/graduate/research-degrees/applying-graduate-research-study-ucl"
$ws.Range("M2").Value = "Entry requirements:/prospective-students/graduate/teacher-training/entry-requirements"
$ws.Range("N2").Value = "How to apply:/prospective-students/graduate/teacher-training/how-apply"
$ws.Range("O2").ClearContents()
